$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the cell to be written as text (shared string) instead of a number,
    # then clear the number-format style so no style index is left on the cell.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

function Set-TextIfNotSet($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Column A (Order Number) - rows 2-7 switch from numeric to text; row 8 unchanged;
# row 9 is new.
Set-TextCell "A2" "1"
Set-TextCell "A3" "2"
Set-TextCell "A4" "3"
Set-TextCell "A5" "4"
Set-TextCell "A6" "5"
Set-TextCell "A7" "7"
Set-TextCell "A9" "9"

# Column B (Part Number) - rows 2-7 switch from numeric to text; row 8 unchanged;
# row 9 is new.
Set-TextCell "B2" "23498345933"
Set-TextCell "B3" "23942039482"
Set-TextCell "B4" "29384023948"
Set-TextCell "B5" "23029849023"
Set-TextCell "B6" "349584398539"
Set-TextCell "B7" "1290138230948"
Set-TextCell "B9" "230924802"

# Column C (Order Details) - rows 2-8 unchanged; row 9 is new.
Set-TextIfNotSet "C9" "cactus"

# Column D (Price) - rows 2,3,4,5,7 switch from numeric to text; rows 6,8 unchanged;
# row 9 is new.
Set-TextCell "D2" "600"
Set-TextCell "D3" "50"
Set-TextCell "D4" "100"
Set-TextCell "D5" "25"
Set-TextCell "D7" "500"
Set-TextCell "D9" "`$300"

# Column E (Seller) - rows 2-8 unchanged; row 9 is new.
Set-TextIfNotSet "E9" "Tim"

# Column F (Buyer) - rows 2-8 unchanged; row 9 is new.
Set-TextIfNotSet "F9" "Mike"
